# Update "想去人数" (F column) values across sheets, matching the
# gh-pages regenerated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 123
$ws1.Range("F4").Value  = 2001
$ws1.Range("F5").Value  = 315
$ws1.Range("F6").Value  = 60
$ws1.Range("F8").Value  = 2045
$ws1.Range("F9").Value  = 10373
$ws1.Range("F11").Value = 150
$ws1.Range("F12").Value = 267
$ws1.Range("F13").Value = 199
$ws1.Range("F15").Value = 7228
$ws1.Range("F17").Value = 685
$ws1.Range("F18").Value = 132
$ws1.Range("F20").Value = 266

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 18

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 123
$ws4.Range("F4").Value  = 2001
$ws4.Range("F5").Value  = 315
$ws4.Range("F6").Value  = 60
$ws4.Range("F7").Value  = 18
$ws4.Range("F9").Value  = 2045
$ws4.Range("F12").Value = 10374
$ws4.Range("F14").Value = 150
$ws4.Range("F15").Value = 267
$ws4.Range("F16").Value = 199
$ws4.Range("F18").Value = 7228
$ws4.Range("F20").Value = 685
$ws4.Range("F21").Value = 132
$ws4.Range("F23").Value = 266
